$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append 5 new rows (157-161) to the device-history master table:
#   regcntr_id=10002, device_id=3000176..3000180
#   lang_code="eng", is_active=TRUE, cr_by="superadmin",
#   cr_dtimes="now()", eff_dtimes="now()"
$startRow = 157
$startDeviceId = 3000176
$rowCount = 5

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = 10002
    $ws.Cells.Item($row, 2).Value = ($startDeviceId + $i)
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $ws.Cells.Item($row, 7).Value = "now()"
}

# Match the author's final view state: scrolled down, B157 selected.
$excel.ActiveWindow.ScrollRow = 152
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B157").Select()
